# Fruta / hortaliza, semanal
# Insert a new weekly record at row 319 (Mango, Vega Central Mapocho de Santiago),
# shifting the existing rows 319:450 down to 320:451.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 319 - pushes old row 319 (and everything below it) down by one.
$ws.Rows.Item(319).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(319, 1).Value  = 9
$ws.Cells.Item(319, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(319, 3).Value  = "Metropolitana"
$ws.Cells.Item(319, 4).Value  = 44755
$ws.Cells.Item(319, 5).Value  = 13
$ws.Cells.Item(319, 6).Value  = "Fruta"
$ws.Cells.Item(319, 7).Value  = 100108
$ws.Cells.Item(319, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(319, 9).Value  = 100108002
$ws.Cells.Item(319, 10).Value = "Mango"
$ws.Cells.Item(319, 11).Value = "Sin especificar"
$ws.Cells.Item(319, 12).Value = "Primera"
$ws.Cells.Item(319, 13).Value = 550
$ws.Cells.Item(319, 14).Value = 7000
$ws.Cells.Item(319, 15).Value = 7500
$ws.Cells.Item(319, 16).Value = 7227
$ws.Cells.Item(319, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(319, 18).Value = "Brasil"
$ws.Cells.Item(319, 19).Value = 1807
$ws.Cells.Item(319, 20).Value = 4
